$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-10-31 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-11-01 Saturday", 2)

# Update the 25 division problems laid out in a 5x5 grid inside the first table.
# The table has 20 rows total; rows 1, 5, 9, 13, 17 contain the actual data
# (the rows in between are blank rows left for handwritten practice).
$t = $d.Tables.Item(1)

$values = @(
    @("39÷9=4, 3", "13÷5=2, 3", "91÷2=45, 1", "49÷2=24, 1", "23÷8=2, 7"),
    @("85÷9=9, 4", "47÷2=23, 1", "31÷9=3, 4", "41÷2=20, 1", "14÷8=1, 6"),
    @("63÷9=7, 0", "69÷9=7, 6", "89÷4=22, 1", "94÷5=18, 4", "21÷3=7, 0"),
    @("15÷2=7, 1", "68÷2=34, 0", "53÷7=7, 4", "75÷2=37, 1", "85÷4=21, 1"),
    @("62÷6=10, 2", "97÷2=48, 1", "85÷7=12, 1", "68÷7=9, 5", "59÷3=19, 2")
)

$dataRows = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $row = $dataRows[$i]
    for ($col = 1; $col -le 5; $col++) {
        $t.Cell($row, $col).Range.Text = $values[$i][$col - 1]
    }
}
